# Re-process the data with the newly curated dimensions.
# The "grupo", "sexo", "descripcion-ocupacion" and "direccion-provincial-nombre"
# columns move from being iaest-dimension (dim) values to iaest-measure (medida)
# values, and their corresponding "mapping-*.xlsx" / URI metadata rows are
# cleared out where no longer applicable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G - grupo
$ws.Cells.Item(2, 7).Value = "iaest-measure:grupo"
$ws.Cells.Item(3, 7).Value = "medida"
$ws.Cells.Item(4, 7).Value = "xsd:int"
$ws.Cells.Item(5, 7).Clear()

# Column H - sexo
$ws.Cells.Item(2, 8).Value = "iaest-measure:sexo"
$ws.Cells.Item(3, 8).Value = "medida"
$ws.Cells.Item(4, 8).Value = "xsd:int"
$ws.Cells.Item(5, 8).Clear()

# Column I - descripcion-ocupacion
$ws.Cells.Item(2, 9).Value = "iaest-measure:descripcion-ocupacion"
$ws.Cells.Item(3, 9).Value = "medida"
$ws.Cells.Item(4, 9).Value = "xsd:int"
$ws.Cells.Item(5, 9).Clear()

# Column K - direccion-provincial-nombre
$ws.Cells.Item(2, 11).Value = "iaest-measure:direccion-provincial-nombre"
$ws.Cells.Item(3, 11).Value = "medida"
$ws.Cells.Item(4, 11).Value = "xsd:int"
